$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Company Name",
    "Company Number",
    "Incorporation Date",
    "Status",
    "Source",
    "Date Downloaded",
    "Time Discovered",
    "Category",
    "SIC Codes",
    "SIC Description",
    "Typical Use Case"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Build the header style on a scratch cell, then stamp it onto the header
# row in a single paste so only one new cellXf combination is recorded
# (bold font + thin box border + center/top alignment), instead of one
# new style per individual property write.
$helper = $ws.Range("Z1")
$helper.Font.Bold = $true
$helper.HorizontalAlignment = -4108
$helper.VerticalAlignment = -4160
$helper.Borders.Weight = 2

$rng = $ws.Range("A1:K1")
$helper.Copy()
$rng.PasteSpecial(-4122)
$helper.Clear()
